# Automatische test-sync: 2025-06-30 19:46:50
# Adds the 5th test-mail row to "Logs" and the matching aggregate row to
# "Dashboard", then widens the dashboard chart series + conditional
# formatting ranges so they keep covering the new data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 6
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A6").Value = "Wat is de levertijd van de EcoPro-700?"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value = "Testmail #5: Wat is de levertijd van de EcoPro-700?"
$logs.Range("D6").Value = "Productinformatie"
$logs.Range("E6").Value = "Beste klant,`nBedankt voor uw vraag over de levertijd van de EcoPro-700. De standaard levertijd voor dit product is momenteel 3-5 werkdagen. `nMocht u verdere vragen hebben of uw bestelling willen opvolgen, laat het ons gerust weten.`nMet vriendelijke groet,`n[Naam Bedrijf]"
$logs.Range("F6").Value = "2025-06-30 19:46:29"
$logs.Range("G6").Value = "Ja"
$logs.Range("H6").Value = "Nee"
$logs.Range("I6").Value = "Ja"
$logs.Range("J6").Value = "Nee"

# Writing a wrapped multi-line value through COM makes the host best-guess
# a custom row height (ht=".." customHeight="1"). Re-autofitting drops the
# custom-height pin again, matching how the original rows (written without
# ever taking this code path) look — no explicit row height at all.
$logs.Rows(6).AutoFit()

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append row 5 (aggregate for the new category)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Productinformatie"
$dash.Range("B5").Value = 1

# ---------------------------------------------------------------------
# 3. Widen the chart series ranges on the dashboard chart so they include
#    the newly added row 5.
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$5,Dashboard!`$B`$2:`$B`$5,1)"

# ---------------------------------------------------------------------
# 4. Extend the conditional-formatting ranges on the Logs sheet from
#    row 5 to row 6 so the new row is covered too.
# ---------------------------------------------------------------------
function Extend-ConditionalFormatting($sheet, $oldRange, $newRange) {
    $fcs = $sheet.Range($oldRange).FormatConditions
    $fc = $fcs.Item(1)
    $fc.ModifyAppliesToRange($sheet.Range($newRange))
}

Extend-ConditionalFormatting $logs "D2:D5" "D2:D6"
Extend-ConditionalFormatting $logs "G2:G5" "G2:G6"
Extend-ConditionalFormatting $logs "H2:H5" "H2:H6"
Extend-ConditionalFormatting $logs "I2:I5" "I2:I6"
Extend-ConditionalFormatting $logs "J2:J5" "J2:J6"
